$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.038.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "'2.508.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'518.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'131.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "'2.507.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'0.0968"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "'5.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("D13").Value = "'0.329"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").Value = "'2.946.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'58.077.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "'21.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "'2.504.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'320.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.69%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'64.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "'7.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "'0.0₃0744"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "'167.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'6.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'17.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").Value = "'36.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "'0.765"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("D42").Value = "'274.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "'3.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "'4.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'0.0917"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'120.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("D49").Value = "'17.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "'16.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.34%  "
